$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2, C2: drop their explicit cell styles (hyperlink style / default) while
# writing the new profile's URL and user id as plain text.
$ws.Range("B2").Clear()
$ws.Range("B2").Value = "http://10.0.74.5/index/home"

$ws.Range("C2").Clear()
$ws.Range("C2").Value = "4fku01"

# D2: drop its "@" text style but keep the value stored as text (not a
# number) by round-tripping the number format.
$ws.Range("D2").Clear()
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "123"
$ws.Range("D2").NumberFormat = "General"

# E2 / G2: same date string repeated, forced to stay text instead of being
# parsed into a date serial.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10/10/2019"
$ws.Range("E2").NumberFormat = "General"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "10/10/2019"
$ws.Range("G2").NumberFormat = "General"

# F2 / H2: new timestamps (already stay text on their own).
$ws.Range("F2").Value = "14:13:10.969"
$ws.Range("H2").Value = "14:13:11.171"

# I2 / J2: new OS / browser strings.
$ws.Range("I2").Value = "Windows Server 2016"
$ws.Range("J2").Value = "Firefox 69.0.2"
